# [LC-850] Release of LetsCo OS v1.3.0
# Rename the generic GPx / BPx KPI identifiers to their zero-padded
# two-digit equivalents (GP1->GP01, GP2->GP02, GP3->GP03,
# BP1->BP01, BP2->BP02, BP3->BP03) and reset the sheet's saved
# view/selection back to cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename GPx / BPx labels in column B -----------------------------
$map = @{
    "GP1" = "GP01"
    "GP2" = "GP02"
    "GP3" = "GP03"
    "BP1" = "BP01"
    "BP2" = "BP02"
    "BP3" = "BP03"
}

$lastRow = $ws.Cells(1, 1).SpecialCells(11).Row
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Text
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}

# --- Reset the view back to A1 ----------------------------------------
$ws.Range("A1").Select()
